# Atualização de bases das ligas, do dia: 12-06-2024 às 23:38
#
# The source feed re-emitted three pairs of fixtures (same Div/Date) with
# their row order flipped. Swap the data payload of each row pair while
# leaving column A (the positional index) untouched. Only the columns
# whose value actually differs between the two rows of a pair are
# touched, so cells that are identical in both rows are left byte-exact.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: the two row numbers, and the 1-based column numbers that
# need to be swapped between them (columns identical in both rows are
# skipped on purpose).
$rowPairs = @(
    @{ Row1 = 102; Row2 = 103; Cols = @(2,5,6,7,8,9,10,12,13,14,15,17,18,19,20,22,23,24,27,29,30) },
    @{ Row1 = 114; Row2 = 115; Cols = @(2,5,6,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,25,27,28,29,30) },
    @{ Row1 = 162; Row2 = 163; Cols = @(2,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,22,23,24,25,27,28) }
)

foreach ($pair in $rowPairs) {
    $r1 = $pair.Row1
    $r2 = $pair.Row2

    foreach ($c in $pair.Cols) {
        $cell1 = $ws.Cells.Item($r1, $c)
        $cell2 = $ws.Cells.Item($r2, $c)

        $v1 = $cell1.Value2
        $v2 = $cell2.Value2

        $cell1.Value2 = $v2
        $cell2.Value2 = $v1
    }
}
